$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 489) holds the "Förändrad" (last changed) date,
# currently stored as serial date 45204 (2023-10-05). Bump it by one day to
# 45205 (2023-10-06) for every data row.
$ws.Range("C2:C489").Value = 45205
